$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.754.94"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.124.70"
$ws.Range("E3").Value = "  +10.60%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "2.429.05"
$ws.Range("E13").Value = "  +10.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").Value = "2.119.79"
$ws.Range("E16").Value = "  +10.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "36.736.28"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "0.0₃0843"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "242.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -7.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.65%  "
$ws.Range("E29").Value = "  -8.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +56.41%  "
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0956"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +17.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.945"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.23%  "
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  -7.68%  "
$ws.Range("E41").Value = "  +7.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("D46").Value = "1.359.30"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("E48").Value = "  +10.55%  "
$ws.Range("D49").Value = "2.317.91"
$ws.Range("E49").Value = "  +10.34%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("E51").Value = "  +0.58%  "
